$d = $word.ActiveDocument

# Locate the paragraph that holds the lone "_GoBack" bookmark: it's the
# empty "List Paragraph"-styled paragraph that follows the finger list
# items (the only empty ListParagraph paragraph in the doc).
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "List Paragraph" -and $p.Range.Text.Trim() -eq "") {
        $target = $p
        break
    }
}

# Drop the explicit "ListParagraph" formatting so the paragraph reverts to
# plain/default (no <w:pPr> at all), matching the rest of the body text.
$target.Format.Style = "Normal"

# Append the new sentence after the bookmark markers (InsertAfter() appends
# at the end of the paragraph's range, i.e. after the existing
# bookmarkStart/bookmarkEnd, keeping their original order).
$tRange = $target.Range
$insertion = $d.Range($tRange.End - 1, $tRange.End - 1)
$insertion.InsertAfter("The thumb is counted every 8 times so I divided 8 into every problem to get my answers.")

# Insert a brand-new empty paragraph right after the one we just edited
# (and before the following blank paragraph), by splitting the paragraph
# mark just before it.
$tRange2 = $target.Range
$splitPoint = $d.Range($tRange2.End - 1, $tRange2.End - 1)
$splitPoint.Text = "`r"
